$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 946.25  # H17: 1007.8205 -> 946.25
$ws.Cells.Item(17, 10).Value = 946.25  # J17: 1007.8205 -> 946.25
$ws.Cells.Item(17, 12).Value = 2838.75  # L17: 3023.4615 -> 2838.75
$ws.Cells.Item(17, 14).Value = -3174.75  # N17: -3359.4615 -> -3174.75
$ws.Cells.Item(33, 8).Value = 495.29166  # H33: 559.8095 -> 495.29166
$ws.Cells.Item(33, 9).Value = 93.411766  # I33: 104.07143 -> 93.411766
$ws.Cells.Item(33, 11).Value = 93.411766  # K33: 104.07143 -> 93.411766
$ws.Cells.Item(33, 13).Value = 135.588234  # M33: 124.92857 -> 135.588234
$ws.Cells.Item(97, 8).Value = 1570.6666  # H97: 1409.75 -> 1570.6666
$ws.Cells.Item(97, 10).Value = 1622.5454  # J97: 1447 -> 1622.5454
$ws.Cells.Item(97, 12).Value = 4867.6362  # L97: 4341 -> 4867.6362
$ws.Cells.Item(97, 14).Value = -5859.6362  # N97: -5333 -> -5859.6362
$ws.Cells.Item(129, 8).Value = 526.5714  # H129: 642.125 -> 526.5714
$ws.Cells.Item(129, 10).Value = 1100  # J129: 1217 -> 1100
$ws.Cells.Item(129, 12).Value = 3300  # L129: 3651 -> 3300
$ws.Cells.Item(129, 14).Value = -13300  # N129: -13651 -> -13300
$ws.Cells.Item(132, 8).Value = 599308.8  # H132: 702028.8 -> 599308.8
$ws.Cells.Item(132, 9).Value = 1734.2712  # I132: 2050.9795 -> 1734.2712
$ws.Cells.Item(132, 10).Value = 2132217.5  # J132: 2335310.5 -> 2132217.5
$ws.Cells.Item(132, 11).Value = 5202.813599999999  # K132: 6152.9385 -> 5202.813599999999
$ws.Cells.Item(132, 12).Value = 6396652.5  # L132: 7005931.5 -> 6396652.5
$ws.Cells.Item(132, 13).Value = -2672.813599999999  # M132: -3622.9385 -> -2672.813599999999
$ws.Cells.Item(132, 14).Value = -6401712.5  # N132: -7010991.5 -> -6401712.5
$ws.Cells.Item(135, 8).Value = 20011  # H135: 20778.154 -> 20011
$ws.Cells.Item(135, 9).Value = 21745.25  # I135: 22687.87 -> 21745.25
$ws.Cells.Item(135, 11).Value = 195707.25  # K135: 204190.83 -> 195707.25
$ws.Cells.Item(135, 13).Value = -193172.25  # M135: -201655.83 -> -193172.25
$ws.Cells.Item(137, 8).Value = 3335812  # H137: 1853418.2 -> 3335812
$ws.Cells.Item(137, 9).Value = 6251672.5  # I137: 2381949.8 -> 6251672.5
$ws.Cells.Item(137, 10).Value = 3400  # J137: 3558.0833 -> 3400
$ws.Cells.Item(137, 11).Value = 18755017.5  # K137: 7145849.399999999 -> 18755017.5
$ws.Cells.Item(137, 12).Value = 10200  # L137: 10674.2499 -> 10200
$ws.Cells.Item(137, 13).Value = -18752467.5  # M137: -7143299.399999999 -> -18752467.5
$ws.Cells.Item(137, 14).Value = -15300  # N137: -15774.2499 -> -15300
$ws.Cells.Item(141, 8).Value = 2978.054  # H141: 3266.5 -> 2978.054
$ws.Cells.Item(141, 9).Value = 1927.7241  # I141: 2048.4285 -> 1927.7241
$ws.Cells.Item(141, 10).Value = 6785.5  # J141: 8950.833000000001 -> 6785.5
$ws.Cells.Item(141, 11).Value = 5783.1723  # K141: 6145.2855 -> 5783.1723
$ws.Cells.Item(141, 12).Value = 20356.5  # L141: 26852.499 -> 20356.5
$ws.Cells.Item(141, 13).Value = -603.1723000000002  # M141: -965.2855 -> -603.1723000000002
$ws.Cells.Item(141, 14).Value = -30716.5  # N141: -37212.499 -> -30716.5

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 2961.51  # H32: 6974.16 -> 2961.51
$ws.Cells.Item(32, 9).Value = 2491.4946  # I32: 3819.6836 -> 2491.4946
$ws.Cells.Item(32, 10).Value = 9206  # J32: 18841 -> 9206
$ws.Cells.Item(32, 11).Value = 2491.4946  # K32: 3819.6836 -> 2491.4946
$ws.Cells.Item(32, 12).Value = 9206  # L32: 18841 -> 9206
$ws.Cells.Item(32, 13).Value = -2204.4946  # M32: -3532.6836 -> -2204.4946
$ws.Cells.Item(32, 14).Value = -9780  # N32: -19415 -> -9780
$ws.Cells.Item(61, 8).Value = 19270654  # H61: 13541580 -> 19270654
$ws.Cells.Item(61, 9).Value = 24415814  # I61: 16146050 -> 24415814
$ws.Cells.Item(61, 10).Value = 93239.27  # J61: 85150.664 -> 93239.27
$ws.Cells.Item(61, 11).Value = 24415814  # K61: 16146050 -> 24415814
$ws.Cells.Item(61, 12).Value = 93239.27  # L61: 85150.664 -> 93239.27
$ws.Cells.Item(61, 13).Value = -24415602  # M61: -16145838 -> -24415602
$ws.Cells.Item(61, 14).Value = -93663.27  # N61: -85574.664 -> -93663.27
$ws.Cells.Item(74, 8).Value = 9316309  # H74: 7453052 -> 9316309
$ws.Cells.Item(74, 9).Value = 13374637  # I74: 9834353 -> 13374637
$ws.Cells.Item(74, 10).Value = 92836.27  # J74: 92669.17999999999 -> 92836.27
$ws.Cells.Item(74, 11).Value = 13374637  # K74: 9834353 -> 13374637
$ws.Cells.Item(74, 12).Value = 92836.27  # L74: 92669.17999999999 -> 92836.27
$ws.Cells.Item(74, 13).Value = -13373763  # M74: -9833479 -> -13373763
$ws.Cells.Item(74, 14).Value = -94584.27  # N74: -94417.17999999999 -> -94584.27
$ws.Cells.Item(77, 8).Value = 9316309  # H77: 7453052 -> 9316309
$ws.Cells.Item(77, 9).Value = 13374637  # I77: 9834353 -> 13374637
$ws.Cells.Item(77, 10).Value = 92836.27  # J77: 92669.17999999999 -> 92836.27
$ws.Cells.Item(77, 11).Value = 66873185  # K77: 49171765 -> 66873185
$ws.Cells.Item(77, 12).Value = 464181.35  # L77: 463345.9 -> 464181.35
$ws.Cells.Item(77, 13).Value = -66868817  # M77: -49167397 -> -66868817
$ws.Cells.Item(77, 14).Value = -472917.35  # N77: -472081.9 -> -472917.35
$ws.Cells.Item(136, 8).Value = 19270654  # H136: 13541580 -> 19270654
$ws.Cells.Item(136, 9).Value = 24415814  # I136: 16146050 -> 24415814
$ws.Cells.Item(136, 10).Value = 93239.27  # J136: 85150.664 -> 93239.27
$ws.Cells.Item(136, 11).Value = 73247442  # K136: 48438150 -> 73247442
$ws.Cells.Item(136, 12).Value = 279717.81  # L136: 255451.992 -> 279717.81
$ws.Cells.Item(136, 13).Value = -73244892  # M136: -48435600 -> -73244892
$ws.Cells.Item(136, 14).Value = -284817.81  # N136: -260551.992 -> -284817.81

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(26, 8).Value = 14924.429  # H26: 11037.1 -> 14924.429
$ws.Cells.Item(26, 9).Value = 14924.429  # I26: 11037.1 -> 14924.429
$ws.Cells.Item(26, 11).Value = 14924.429  # K26: 11037.1 -> 14924.429
$ws.Cells.Item(26, 13).Value = -14632.429  # M26: -10745.1 -> -14632.429
$ws.Cells.Item(96, 8).Value = 25000  # H96: 8920.727999999999 -> 25000
$ws.Cells.Item(96, 9).Value = 25000  # I96: 2916 -> 25000
$ws.Cells.Item(96, 10).Value = 25000  # J96: 24933.334 -> 25000
$ws.Cells.Item(96, 11).Value = 25000  # K96: 2916 -> 25000
$ws.Cells.Item(96, 12).Value = 25000  # L96: 24933.334 -> 25000
$ws.Cells.Item(96, 13).Value = -22254  # M96: -170 -> -22254
$ws.Cells.Item(96, 14).Value = -30492  # N96: -30425.334 -> -30492

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(62, 8).Value = 2475  # H62: 2613.9285 -> 2475
$ws.Cells.Item(62, 9).Value = 2376.25  # I62: 2833.3333 -> 2376.25
$ws.Cells.Item(62, 10).Value = 2507.9167  # J62: 2554.0908 -> 2507.9167
$ws.Cells.Item(62, 11).Value = 2376.25  # K62: 2833.3333 -> 2376.25
$ws.Cells.Item(62, 12).Value = 2507.9167  # L62: 2554.0908 -> 2507.9167
$ws.Cells.Item(62, 13).Value = -1752.25  # M62: -2209.3333 -> -1752.25
$ws.Cells.Item(62, 14).Value = -3755.9167  # N62: -3802.0908 -> -3755.9167
$ws.Cells.Item(65, 8).Value = 2475  # H65: 2613.9285 -> 2475
$ws.Cells.Item(65, 9).Value = 2376.25  # I65: 2833.3333 -> 2376.25
$ws.Cells.Item(65, 10).Value = 2507.9167  # J65: 2554.0908 -> 2507.9167
$ws.Cells.Item(65, 11).Value = 11881.25  # K65: 14166.6665 -> 11881.25
$ws.Cells.Item(65, 12).Value = 12539.5835  # L65: 12770.454 -> 12539.5835
$ws.Cells.Item(65, 13).Value = -8761.25  # M65: -11046.6665 -> -8761.25
$ws.Cells.Item(65, 14).Value = -18779.5835  # N65: -19010.454 -> -18779.5835
$ws.Cells.Item(132, 8).Value = 15277.085  # H132: 19379.715 -> 15277.085
$ws.Cells.Item(132, 9).Value = 1045.8392  # I132: 1354.4286 -> 1045.8392
$ws.Cells.Item(132, 10).Value = 68407.07000000001  # J132: 73455.57000000001 -> 68407.07000000001
$ws.Cells.Item(132, 11).Value = 3137.5176  # K132: 4063.2858 -> 3137.5176
$ws.Cells.Item(132, 12).Value = 205221.21  # L132: 220366.71 -> 205221.21
$ws.Cells.Item(132, 13).Value = -607.5175999999997  # M132: -1533.2858 -> -607.5175999999997
$ws.Cells.Item(132, 14).Value = -210281.21  # N132: -225426.71 -> -210281.21
$ws.Cells.Item(134, 8).Value = 15913.027  # H134: 17316.373 -> 15913.027
$ws.Cells.Item(134, 9).Value = 1095.2363  # I134: 1185.68 -> 1095.2363
$ws.Cells.Item(134, 10).Value = 61189.61  # J134: 64759.59 -> 61189.61
$ws.Cells.Item(134, 11).Value = 3285.7089  # K134: 3557.04 -> 3285.7089
$ws.Cells.Item(134, 12).Value = 183568.83  # L134: 194278.77 -> 183568.83
$ws.Cells.Item(134, 13).Value = -750.7089000000001  # M134: -1022.04 -> -750.7089000000001
$ws.Cells.Item(134, 14).Value = -188638.83  # N134: -199348.77 -> -188638.83

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(3, 8).Value = 7644.2856  # H3: 5924.2856 -> 7644.2856
$ws.Cells.Item(3, 9).Value = 5585  # I3: 5303.077 -> 5585
$ws.Cells.Item(3, 10).Value = 20000  # J3: 14000 -> 20000
$ws.Cells.Item(3, 11).Value = 16755  # K3: 15909.231 -> 16755
$ws.Cells.Item(3, 12).Value = 60000  # L3: 42000 -> 60000
$ws.Cells.Item(3, 13).Value = -16643  # M3: -15797.231 -> -16643
$ws.Cells.Item(3, 14).Value = -60224  # N3: -42224 -> -60224
$ws.Cells.Item(5, 8).Value = 821  # H5: 779.4 -> 821
$ws.Cells.Item(5, 10).Value = 1247.5  # J5: 966.6667 -> 1247.5
$ws.Cells.Item(5, 12).Value = 3742.5  # L5: 2900.0001 -> 3742.5
$ws.Cells.Item(5, 14).Value = -3966.5  # N5: -3124.0001 -> -3966.5
$ws.Cells.Item(97, 8).Value = 1841  # H97: 1844 -> 1841
$ws.Cells.Item(97, 9).Value = 3411.1428  # I97: 3009.75 -> 3411.1428
$ws.Cells.Item(97, 10).Value = 467.125  # J97: 511.7143 -> 467.125
$ws.Cells.Item(97, 11).Value = 10233.4284  # K97: 9029.25 -> 10233.4284
$ws.Cells.Item(97, 12).Value = 1401.375  # L97: 1535.1429 -> 1401.375
$ws.Cells.Item(97, 13).Value = -9737.428400000001  # M97: -8533.25 -> -9737.428400000001
$ws.Cells.Item(97, 14).Value = -2393.375  # N97: -2527.1429 -> -2393.375
$ws.Cells.Item(113, 8).Value = 418.93332  # H113: 451.17392 -> 418.93332
$ws.Cells.Item(113, 9).Value = 328.94736  # I113: 338.25 -> 328.94736
$ws.Cells.Item(113, 11).Value = 986.84208  # K113: 1014.75 -> 986.84208
$ws.Cells.Item(113, 13).Value = 1183.15792  # M113: 1155.25 -> 1183.15792
$ws.Cells.Item(121, 8).Value = 39236960  # H121: 48529704 -> 39236960
$ws.Cells.Item(121, 9).Value = 1471.6666  # I121: 1832.5 -> 1471.6666
$ws.Cells.Item(121, 10).Value = 44978740  # J121: 54238868 -> 44978740
$ws.Cells.Item(121, 11).Value = 4414.9998  # K121: 5497.5 -> 4414.9998
$ws.Cells.Item(121, 12).Value = 134936220  # L121: 162716604 -> 134936220
$ws.Cells.Item(121, 13).Value = -3104.9998  # M121: -4187.5 -> -3104.9998
$ws.Cells.Item(121, 14).Value = -134938840  # N121: -162719224 -> -134938840
$ws.Cells.Item(135, 8).Value = 821  # H135: 779.4 -> 821
$ws.Cells.Item(135, 10).Value = 1247.5  # J135: 966.6667 -> 1247.5
$ws.Cells.Item(135, 12).Value = 11227.5  # L135: 8700.0003 -> 11227.5
$ws.Cells.Item(135, 14).Value = -16297.5  # N135: -13770.0003 -> -16297.5

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 23948.133  # H132: 23424.39 -> 23948.133
$ws.Cells.Item(132, 9).Value = 1325  # I132: 1248.7812 -> 1325
$ws.Cells.Item(132, 10).Value = 64952.562  # J132: 74111.5 -> 64952.562
$ws.Cells.Item(132, 11).Value = 3975  # K132: 3746.3436 -> 3975
$ws.Cells.Item(132, 12).Value = 194857.686  # L132: 222334.5 -> 194857.686
$ws.Cells.Item(132, 13).Value = -1445  # M132: -1216.3436 -> -1445
$ws.Cells.Item(132, 14).Value = -199917.686  # N132: -227394.5 -> -199917.686

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 23177.674  # H132: 26110.152 -> 23177.674
$ws.Cells.Item(132, 9).Value = 14991.739  # I132: 16949.902 -> 14991.739
$ws.Cells.Item(132, 10).Value = 51419.15  # J132: 57153.223 -> 51419.15
$ws.Cells.Item(132, 11).Value = 44975.217  # K132: 50849.70599999999 -> 44975.217
$ws.Cells.Item(132, 12).Value = 154257.45  # L132: 171459.669 -> 154257.45
$ws.Cells.Item(132, 13).Value = -42445.217  # M132: -48319.70599999999 -> -42445.217
$ws.Cells.Item(132, 14).Value = -159317.45  # N132: -176519.669 -> -159317.45
$ws.Cells.Item(136, 8).Value = 30485.479  # H136: 38937.465 -> 30485.479
$ws.Cells.Item(136, 9).Value = 21487.5  # I136: 29437.656 -> 21487.5
$ws.Cells.Item(136, 10).Value = 51052.285  # J136: 56437.105 -> 51052.285
$ws.Cells.Item(136, 11).Value = 64462.5  # K136: 88312.96799999999 -> 64462.5
$ws.Cells.Item(136, 12).Value = 153156.855  # L136: 169311.315 -> 153156.855
$ws.Cells.Item(136, 13).Value = -61912.5  # M136: -85762.96799999999 -> -61912.5
$ws.Cells.Item(136, 14).Value = -158256.855  # N136: -174411.315 -> -158256.855
$ws.Cells.Item(138, 8).Value = 0  # H138: 44041 -> 0
$ws.Cells.Item(138, 10).Value = 0  # J138: 44041 -> 0
$ws.Cells.Item(138, 12).Value = 0  # L138: 44041 -> 0
$ws.Cells.Item(138, 14).ClearContents()  # N138: remove (was -54321)
